$d = $word.ActiveDocument
$rng = $d.Content
$found = $rng.Find.Execute("For the missing-data dataset XXX.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target paragraph text"
}
$s = $rng.Start
$rng.Text = ""
$insertRng = $d.Range($s, $s)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">For the missing-data dataset 250 post-warmup samples were taken from 4</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">indpendent Markov chains after 100 warmup samples. The sampling was initialised</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">using the mass matric from the complete measurement dataset and the warmup</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">consisted of step size adaption for 100 samples. The resulting posterior distribution</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">had an</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:acc><m:accPr><m:chr m:val="̂"/></m:accPr><m:e><m:r><m:t>R</m:t></m:r></m:e></m:acc><m:r><m:rPr><m:sty m:val="p"/></m:rPr><m:t>=</m:t></m:r><m:r><m:t>1.01</m:t></m:r></m:oMath><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">for the log-probability and did not exhibit post-warmup</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">divergences that were not a result of differential equation errors.</w:t></w:r></w:p>'
$insertRng.InsertXML($xml)
Write-Host "Paragraph updated"
